# ============================================================================
# Scheduled market-data refresh for the Halicarnassus Leve-profit workbook.
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) on each
# job-sheet with freshly polled Universalis price data, row by row.
# ============================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: ALC
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")
# row 9: H9=1315, I9=630, J9=2000, K9=630, L9=2000, M9=-461, N9=-2338
$ws.Cells.Item(9, 8).Value = 1315
$ws.Cells.Item(9, 9).Value = 630
$ws.Cells.Item(9, 10).Value = 2000
$ws.Cells.Item(9, 11).Value = 630
$ws.Cells.Item(9, 12).Value = 2000
$ws.Cells.Item(9, 13).Value = -461
$ws.Cells.Item(9, 14).Value = -2338
# row 20: H20=1750, I20=1750, K20=1750, M20=-1520
$ws.Cells.Item(20, 8).Value = 1750
$ws.Cells.Item(20, 9).Value = 1750
$ws.Cells.Item(20, 11).Value = 1750
$ws.Cells.Item(20, 13).Value = -1520
# row 35: H35=1750, I35=1750, K35=1750, M35=-1371
$ws.Cells.Item(35, 8).Value = 1750
$ws.Cells.Item(35, 9).Value = 1750
$ws.Cells.Item(35, 11).Value = 1750
$ws.Cells.Item(35, 13).Value = -1371
# row 62: H62=20000, I62=0, J62=20000, K62=0, L62=20000, M62=(clear), N62=-21248
$ws.Cells.Item(62, 8).Value = 20000
$ws.Cells.Item(62, 9).Value = 0
$ws.Cells.Item(62, 10).Value = 20000
$ws.Cells.Item(62, 11).Value = 0
$ws.Cells.Item(62, 12).Value = 20000
$ws.Cells.Item(62, 13).ClearContents()
$ws.Cells.Item(62, 14).Value = -21248
# row 65: H65=20000, I65=0, J65=20000, K65=0, L65=100000, M65=(clear), N65=-106240
$ws.Cells.Item(65, 8).Value = 20000
$ws.Cells.Item(65, 9).Value = 0
$ws.Cells.Item(65, 10).Value = 20000
$ws.Cells.Item(65, 11).Value = 0
$ws.Cells.Item(65, 12).Value = 100000
$ws.Cells.Item(65, 13).ClearContents()
$ws.Cells.Item(65, 14).Value = -106240
# row 74: H74=12467, I74=3678.5, K74=3678.5, M74=-2742.5
$ws.Cells.Item(74, 8).Value = 12467
$ws.Cells.Item(74, 9).Value = 3678.5
$ws.Cells.Item(74, 11).Value = 3678.5
$ws.Cells.Item(74, 13).Value = -2742.5
# row 77: H77=12467, I77=3678.5, K77=18392.5, M77=-13712.5
$ws.Cells.Item(77, 8).Value = 12467
$ws.Cells.Item(77, 9).Value = 3678.5
$ws.Cells.Item(77, 11).Value = 18392.5
$ws.Cells.Item(77, 13).Value = -13712.5
# row 96: H96=418.22223, I96=302, K96=906, M96=467
$ws.Cells.Item(96, 8).Value = 418.22223
$ws.Cells.Item(96, 9).Value = 302
$ws.Cells.Item(96, 11).Value = 906
$ws.Cells.Item(96, 13).Value = 467
# row 116: H116=4353, I116=4126.143, J116=4750, K116=4126.143, L116=4750, M116=-684.143, N116=-11634
$ws.Cells.Item(116, 8).Value = 4353
$ws.Cells.Item(116, 9).Value = 4126.143
$ws.Cells.Item(116, 10).Value = 4750
$ws.Cells.Item(116, 11).Value = 4126.143
$ws.Cells.Item(116, 12).Value = 4750
$ws.Cells.Item(116, 13).Value = -684.143
$ws.Cells.Item(116, 14).Value = -11634
# row 125: H125=4027.5557, I125=3749.7144, K125=33747.4296, M125=-31287.4296
$ws.Cells.Item(125, 8).Value = 4027.5557
$ws.Cells.Item(125, 9).Value = 3749.7144
$ws.Cells.Item(125, 11).Value = 33747.4296
$ws.Cells.Item(125, 13).Value = -31287.4296
# row 132: H132=22164.75, I132=30096.715, J132=11060, K132=90290.145, L132=33180, M132=-87760.145, N132=-38240
$ws.Cells.Item(132, 8).Value = 22164.75
$ws.Cells.Item(132, 9).Value = 30096.715
$ws.Cells.Item(132, 10).Value = 11060
$ws.Cells.Item(132, 11).Value = 90290.145
$ws.Cells.Item(132, 12).Value = 33180
$ws.Cells.Item(132, 13).Value = -87760.145
$ws.Cells.Item(132, 14).Value = -38240
# row 133: H133=127186.664, J133=127186.664, L133=127186.664, N133=-137306.664
$ws.Cells.Item(133, 8).Value = 127186.664
$ws.Cells.Item(133, 10).Value = 127186.664
$ws.Cells.Item(133, 12).Value = 127186.664
$ws.Cells.Item(133, 14).Value = -137306.664

# ---------------------------------------------------------------------------
# Sheet: ARM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")
# row 2: H2=2366.9412, I2=2554, J2=1918, K2=2554, L2=1918, M2=-2441, N2=-2144
$ws.Cells.Item(2, 8).Value = 2366.9412
$ws.Cells.Item(2, 9).Value = 2554
$ws.Cells.Item(2, 10).Value = 1918
$ws.Cells.Item(2, 11).Value = 2554
$ws.Cells.Item(2, 12).Value = 1918
$ws.Cells.Item(2, 13).Value = -2441
$ws.Cells.Item(2, 14).Value = -2144
# row 18: H18=3969, J18=3969, L18=3969, N18=-4613
$ws.Cells.Item(18, 8).Value = 3969
$ws.Cells.Item(18, 10).Value = 3969
$ws.Cells.Item(18, 12).Value = 3969
$ws.Cells.Item(18, 14).Value = -4613
# row 61: H61=5681.8667, I61=5267.727, K61=5267.727, M61=-5055.727
$ws.Cells.Item(61, 8).Value = 5681.8667
$ws.Cells.Item(61, 9).Value = 5267.727
$ws.Cells.Item(61, 11).Value = 5267.727
$ws.Cells.Item(61, 13).Value = -5055.727
# row 76: H76=0, J76=0, L76=0, N76=(clear)
$ws.Cells.Item(76, 8).Value = 0
$ws.Cells.Item(76, 10).Value = 0
$ws.Cells.Item(76, 12).Value = 0
$ws.Cells.Item(76, 14).ClearContents()
# row 79: H79=0, J79=0, L79=0, N79=(clear)
$ws.Cells.Item(79, 8).Value = 0
$ws.Cells.Item(79, 10).Value = 0
$ws.Cells.Item(79, 12).Value = 0
$ws.Cells.Item(79, 14).ClearContents()
# row 86: H86=0, I86=0, J86=0, K86=0, L86=0, M86=(clear), N86=(clear)
$ws.Cells.Item(86, 8).Value = 0
$ws.Cells.Item(86, 9).Value = 0
$ws.Cells.Item(86, 10).Value = 0
$ws.Cells.Item(86, 11).Value = 0
$ws.Cells.Item(86, 12).Value = 0
$ws.Cells.Item(86, 13).ClearContents()
$ws.Cells.Item(86, 14).ClearContents()
# row 88: H88=2133.1667, J88=2159.8, L88=2159.8, N88=-2971.8
$ws.Cells.Item(88, 8).Value = 2133.1667
$ws.Cells.Item(88, 10).Value = 2159.8
$ws.Cells.Item(88, 12).Value = 2159.8
$ws.Cells.Item(88, 14).Value = -2971.8
# row 89: H89=0, I89=0, J89=0, K89=0, L89=0, M89=(clear), N89=(clear)
$ws.Cells.Item(89, 8).Value = 0
$ws.Cells.Item(89, 9).Value = 0
$ws.Cells.Item(89, 10).Value = 0
$ws.Cells.Item(89, 11).Value = 0
$ws.Cells.Item(89, 12).Value = 0
$ws.Cells.Item(89, 13).ClearContents()
$ws.Cells.Item(89, 14).ClearContents()
# row 91: H91=2133.1667, J91=2159.8, L91=2159.8, N91=-4967.8
$ws.Cells.Item(91, 8).Value = 2133.1667
$ws.Cells.Item(91, 10).Value = 2159.8
$ws.Cells.Item(91, 12).Value = 2159.8
$ws.Cells.Item(91, 14).Value = -4967.8
# row 116: H116=2366.9412, I116=2554, J116=1918, K116=2554, L116=1918, M116=-260, N116=-6506
$ws.Cells.Item(116, 8).Value = 2366.9412
$ws.Cells.Item(116, 9).Value = 2554
$ws.Cells.Item(116, 10).Value = 1918
$ws.Cells.Item(116, 11).Value = 2554
$ws.Cells.Item(116, 12).Value = 1918
$ws.Cells.Item(116, 13).Value = -260
$ws.Cells.Item(116, 14).Value = -6506
# row 117: H117=48000, I117=0, K117=0, M117=(clear)
$ws.Cells.Item(117, 8).Value = 48000
$ws.Cells.Item(117, 9).Value = 0
$ws.Cells.Item(117, 11).Value = 0
$ws.Cells.Item(117, 13).ClearContents()
# row 136: H136=5681.8667, I136=5267.727, K136=15803.181, M136=-13253.181
$ws.Cells.Item(136, 8).Value = 5681.8667
$ws.Cells.Item(136, 9).Value = 5267.727
$ws.Cells.Item(136, 11).Value = 15803.181
$ws.Cells.Item(136, 13).Value = -13253.181
# row 138: H138=0, J138=0, L138=0, N138=(clear)
$ws.Cells.Item(138, 8).Value = 0
$ws.Cells.Item(138, 10).Value = 0
$ws.Cells.Item(138, 12).Value = 0
$ws.Cells.Item(138, 14).ClearContents()

# ---------------------------------------------------------------------------
# Sheet: BSM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")
# row 3: H3=2366.9412, I3=2554, J3=1918, K3=2554, L3=1918, M3=-2440, N3=-2146
$ws.Cells.Item(3, 8).Value = 2366.9412
$ws.Cells.Item(3, 9).Value = 2554
$ws.Cells.Item(3, 10).Value = 1918
$ws.Cells.Item(3, 11).Value = 2554
$ws.Cells.Item(3, 12).Value = 1918
$ws.Cells.Item(3, 13).Value = -2440
$ws.Cells.Item(3, 14).Value = -2146
# row 99: H99=4686.4443, J99=5326.3335, L99=5326.3335, N99=-8322.333500000001
$ws.Cells.Item(99, 8).Value = 4686.4443
$ws.Cells.Item(99, 10).Value = 5326.3335
$ws.Cells.Item(99, 12).Value = 5326.3335
$ws.Cells.Item(99, 14).Value = -8322.333500000001
# row 107: H107=5448.6313, I107=1315.875, J107=8454.272000000001, K107=1315.875, L107=8454.272000000001, M107=604.125, N107=-12294.272
$ws.Cells.Item(107, 8).Value = 5448.6313
$ws.Cells.Item(107, 9).Value = 1315.875
$ws.Cells.Item(107, 10).Value = 8454.272000000001
$ws.Cells.Item(107, 11).Value = 1315.875
$ws.Cells.Item(107, 12).Value = 8454.272000000001
$ws.Cells.Item(107, 13).Value = 604.125
$ws.Cells.Item(107, 14).Value = -12294.272

# ---------------------------------------------------------------------------
# Sheet: CUL
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")
# row 82: H82=10000, I82=0, K82=0, M82=(clear)
$ws.Cells.Item(82, 8).Value = 10000
$ws.Cells.Item(82, 9).Value = 0
$ws.Cells.Item(82, 11).Value = 0
$ws.Cells.Item(82, 13).ClearContents()
# row 85: H85=10000, I85=0, K85=0, M85=(clear)
$ws.Cells.Item(85, 8).Value = 10000
$ws.Cells.Item(85, 9).Value = 0
$ws.Cells.Item(85, 11).Value = 0
$ws.Cells.Item(85, 13).ClearContents()
# row 132: H132=2080.6667, I132=1591.1428, J132=2509, K132=14320.2852, L132=22581, M132=-11790.2852, N132=-27641
$ws.Cells.Item(132, 8).Value = 2080.6667
$ws.Cells.Item(132, 9).Value = 1591.1428
$ws.Cells.Item(132, 10).Value = 2509
$ws.Cells.Item(132, 11).Value = 14320.2852
$ws.Cells.Item(132, 12).Value = 22581
$ws.Cells.Item(132, 13).Value = -11790.2852
$ws.Cells.Item(132, 14).Value = -27641
# row 136: H136=5232, I136=464.5, K136=1393.5, M136=3706.5
$ws.Cells.Item(136, 8).Value = 5232
$ws.Cells.Item(136, 9).Value = 464.5
$ws.Cells.Item(136, 11).Value = 1393.5
$ws.Cells.Item(136, 13).Value = 3706.5
# row 137: H137=1474, I137=665.6667, J137=3899, K137=1997.0001, L137=11697, M137=3102.9999, N137=-21897
$ws.Cells.Item(137, 8).Value = 1474
$ws.Cells.Item(137, 9).Value = 665.6667
$ws.Cells.Item(137, 10).Value = 3899
$ws.Cells.Item(137, 11).Value = 1997.0001
$ws.Cells.Item(137, 12).Value = 11697
$ws.Cells.Item(137, 13).Value = 3102.9999
$ws.Cells.Item(137, 14).Value = -21897
# row 138: H138=8147.6665, I138=4443.3335, K138=13330.0005, M138=-8190.000499999998
$ws.Cells.Item(138, 8).Value = 8147.6665
$ws.Cells.Item(138, 9).Value = 4443.3335
$ws.Cells.Item(138, 11).Value = 13330.0005
$ws.Cells.Item(138, 13).Value = -8190.000499999998

# ---------------------------------------------------------------------------
# Sheet: GSM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")
# row 47: H47=21834, J47=24736, L47=24736, N47=-25872
$ws.Cells.Item(47, 8).Value = 21834
$ws.Cells.Item(47, 10).Value = 24736
$ws.Cells.Item(47, 12).Value = 24736
$ws.Cells.Item(47, 14).Value = -25872
# row 55: H55=5126.857, I55=5399.8, J55=4444.5, K55=5399.8, L55=4444.5, M55=-5072.8, N55=-5098.5
$ws.Cells.Item(55, 8).Value = 5126.857
$ws.Cells.Item(55, 9).Value = 5399.8
$ws.Cells.Item(55, 10).Value = 4444.5
$ws.Cells.Item(55, 11).Value = 5399.8
$ws.Cells.Item(55, 12).Value = 4444.5
$ws.Cells.Item(55, 13).Value = -5072.8
$ws.Cells.Item(55, 14).Value = -5098.5
# row 122: H122=4824.2, I122=4086, J122=7777, K122=12258, L122=23331, M122=-9808, N122=-28231
$ws.Cells.Item(122, 8).Value = 4824.2
$ws.Cells.Item(122, 9).Value = 4086
$ws.Cells.Item(122, 10).Value = 7777
$ws.Cells.Item(122, 11).Value = 12258
$ws.Cells.Item(122, 12).Value = 23331
$ws.Cells.Item(122, 13).Value = -9808
$ws.Cells.Item(122, 14).Value = -28231

# ---------------------------------------------------------------------------
# Sheet: LTW
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")
# row 40: H40=10000, I40=10000, J40=0, K40=10000, L40=0, M40=-9864, N40=(clear)
$ws.Cells.Item(40, 8).Value = 10000
$ws.Cells.Item(40, 9).Value = 10000
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 11).Value = 10000
$ws.Cells.Item(40, 12).Value = 0
$ws.Cells.Item(40, 13).Value = -9864
$ws.Cells.Item(40, 14).ClearContents()
# row 61: H61=6804.154, I61=5794.6665, K61=5794.6665, M61=-5592.6665
$ws.Cells.Item(61, 8).Value = 6804.154
$ws.Cells.Item(61, 9).Value = 5794.6665
$ws.Cells.Item(61, 11).Value = 5794.6665
$ws.Cells.Item(61, 13).Value = -5592.6665
# row 100: H100=6065.2383, I100=2807.625, J100=8069.923, K100=2807.625, L100=8069.923, M100=-2266.625, N100=-9151.922999999999
$ws.Cells.Item(100, 8).Value = 6065.2383
$ws.Cells.Item(100, 9).Value = 2807.625
$ws.Cells.Item(100, 10).Value = 8069.923
$ws.Cells.Item(100, 11).Value = 2807.625
$ws.Cells.Item(100, 12).Value = 8069.923
$ws.Cells.Item(100, 13).Value = -2266.625
$ws.Cells.Item(100, 14).Value = -9151.922999999999
# row 113: H113=6804.154, I113=5794.6665, K113=5794.6665, M113=-3624.6665
$ws.Cells.Item(113, 8).Value = 6804.154
$ws.Cells.Item(113, 9).Value = 5794.6665
$ws.Cells.Item(113, 11).Value = 5794.6665
$ws.Cells.Item(113, 13).Value = -3624.6665
# row 118: H118=76469.664, J118=76469.664, L118=76469.664, N118=-79783.664
$ws.Cells.Item(118, 8).Value = 76469.664
$ws.Cells.Item(118, 10).Value = 76469.664
$ws.Cells.Item(118, 12).Value = 76469.664
$ws.Cells.Item(118, 14).Value = -79783.664
# row 122: H122=2639.6667, I122=2639.6667, K122=7919.000100000001, M122=-5469.000100000001
$ws.Cells.Item(122, 8).Value = 2639.6667
$ws.Cells.Item(122, 9).Value = 2639.6667
$ws.Cells.Item(122, 11).Value = 7919.000100000001
$ws.Cells.Item(122, 13).Value = -5469.000100000001
# row 128: H128=28000, J128=28000, L128=28000, N128=-37960
$ws.Cells.Item(128, 8).Value = 28000
$ws.Cells.Item(128, 10).Value = 28000
$ws.Cells.Item(128, 12).Value = 28000
$ws.Cells.Item(128, 14).Value = -37960
# row 138: H138=89999, J138=89999, L138=89999, N138=-100279
$ws.Cells.Item(138, 8).Value = 89999
$ws.Cells.Item(138, 10).Value = 89999
$ws.Cells.Item(138, 12).Value = 89999
$ws.Cells.Item(138, 14).Value = -100279

# ---------------------------------------------------------------------------
# Sheet: WVR
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")
# row 92: H92=69500, J92=75000, L92=75000, N92=-79992
$ws.Cells.Item(92, 8).Value = 69500
$ws.Cells.Item(92, 10).Value = 75000
$ws.Cells.Item(92, 12).Value = 75000
$ws.Cells.Item(92, 14).Value = -79992
# row 116: H116=29995, J116=29995, L116=29995, N116=-39173
$ws.Cells.Item(116, 8).Value = 29995
$ws.Cells.Item(116, 10).Value = 29995
$ws.Cells.Item(116, 12).Value = 29995
$ws.Cells.Item(116, 14).Value = -39173
# row 122: H122=0, I122=0, K122=0, M122=(clear)
$ws.Cells.Item(122, 8).Value = 0
$ws.Cells.Item(122, 9).Value = 0
$ws.Cells.Item(122, 11).Value = 0
$ws.Cells.Item(122, 13).ClearContents()

